# Purchase order changes - second step
#
# The "REFERENCIA" column (column D) is removed from the volumetry format
# sheet. Deleting the entire column shifts every column from E onward one
# position to the left (so the per-material FÁBRICA/INSTALACIÓN groupings
# that used to live in E:R now live in D:Q), updates the merged header
# cells accordingly, and drops the now-unused "REFERENCIA" shared string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the REFERENCIA column entirely (was column D).
$ws.Range("D1").EntireColumn.Delete() | Out-Null

# Match the author's resulting selection/active cell.
$ws.Range("F11").Select() | Out-Null
